$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-20 Sunday", "2024-10-21 Monday"),
    @("940×6=5640", "780×9=7020"),
    @("915×7=6405", "725×8=5800"),
    @("979×9=8811", "246×6=1476"),
    @("779×5=3895", "631×7=4417"),
    @("904×8=7232", "554×9=4986"),
    @("705×6=4230", "860×2=1720"),
    @("346×7=2422", "136×6=816"),
    @("697×7=4879", "917×9=8253"),
    @("125×2=250", "404×7=2828"),
    @("268×4=1072", "942×8=7536"),
    @("728×9=6552", "886×8=7088"),
    @("486×2=972", "243×2=486"),
    @("977×6=5862", "313×7=2191"),
    @("913×4=3652", "167×8=1336"),
    @("147×2=294", "797×3=2391"),
    @("181×7=1267", "362×9=3258"),
    @("151×2=302", "431×3=1293"),
    @("260×4=1040", "447×5=2235"),
    @("990×4=3960", "499×6=2994"),
    @("893×7=6251", "119×3=357"),
    @("410×2=820", "970×9=8730"),
    @("317×6=1902", "478×4=1912"),
    @("257×8=2056", "679×8=5432"),
    @("521×7=3647", "804×2=1608"),
    @("220×8=1760", "463×3=1389")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "done"
